$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DATA")

# Fix row 4 (test "loginTestUsingExcelData" third data row):
# browser should be "chrome" (was "firefox") and version should be "N/A" (was "95.0.2")
$ws.Range("C4").Value = "chrome"
$ws.Range("D4").Value = "N/A"

# Fix row 7 (test "newTestUsingExcelData" third data row):
# browser should be "firefox" (was "chrome")
$ws.Range("C7").Value = "firefox"

# Remove the "fname" column (column G) entirely, it is no longer used
$ws.Columns.Item(7).Delete()
